$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# The "Data" (GDP per Capita) values are stored as shared strings (text) in the
# source workbook, not numbers. Mark the full E column range as Text before
# writing so the new values serialize the same way (t="s") as the original file.
$ws.Range("E2:E68").NumberFormat = "@"

# Update existing Data (GDP per Capita) values for years 1950-2010 (rows 2-62)
$ws.Cells.Item(2, 5).Value = "1095"
$ws.Cells.Item(3, 5).Value = "1023"
$ws.Cells.Item(4, 5).Value = "1058"
$ws.Cells.Item(5, 5).Value = "1076"
$ws.Cells.Item(6, 5).Value = "1033"
$ws.Cells.Item(7, 5).Value = "1071"
$ws.Cells.Item(8, 5).Value = "1100"
$ws.Cells.Item(9, 5).Value = "1116"
$ws.Cells.Item(10, 5).Value = "1092"
$ws.Cells.Item(11, 5).Value = "1116"
$ws.Cells.Item(12, 5).Value = "1137"
$ws.Cells.Item(13, 5).Value = "1093"
$ws.Cells.Item(14, 5).Value = "1106"
$ws.Cells.Item(15, 5).Value = "1197"
$ws.Cells.Item(16, 5).Value = "1251"
$ws.Cells.Item(17, 5).Value = "1242"
$ws.Cells.Item(18, 5).Value = "1280"
$ws.Cells.Item(19, 5).Value = "1310"
$ws.Cells.Item(20, 5).Value = "1304"
$ws.Cells.Item(21, 5).Value = "1404"
$ws.Cells.Item(22, 5).Value = "1382"
$ws.Cells.Item(23, 5).Value = "1385"
$ws.Cells.Item(24, 5).Value = "1364"
$ws.Cells.Item(25, 5).Value = "1331"
$ws.Cells.Item(26, 5).Value = "1302"
$ws.Cells.Item(27, 5).Value = "1243"
$ws.Cells.Item(28, 5).Value = "1219"
$ws.Cells.Item(29, 5).Value = "1207"
$ws.Cells.Item(30, 5).Value = "1111"
$ws.Cells.Item(31, 5).Value = "966"
$ws.Cells.Item(32, 5).Value = "912"
$ws.Cells.Item(33, 5).Value = "923"
$ws.Cells.Item(34, 5).Value = "972"
$ws.Cells.Item(35, 5).Value = "1014"
$ws.Cells.Item(36, 5).Value = "897"
$ws.Cells.Item(37, 5).Value = "886"
$ws.Cells.Item(38, 5).Value = "858"
$ws.Cells.Item(39, 5).Value = "877"
$ws.Cells.Item(40, 5).Value = "902"
$ws.Cells.Item(41, 5).Value = "929"
$ws.Cells.Item(42, 5).Value = "932"
$ws.Cells.Item(43, 5).Value = "915.4053283443"
$ws.Cells.Item(44, 5).Value = "934.826261105848"
$ws.Cells.Item(45, 5).Value = "959.18814301581"
$ws.Cells.Item(46, 5).Value = "994.533391761838"
$ws.Cells.Item(47, 5).Value = "1054.71392960536"
$ws.Cells.Item(48, 5).Value = "1091.52752640033"
$ws.Cells.Item(49, 5).Value = "1107.50094097432"
$ws.Cells.Item(50, 5).Value = "1137.63815124113"
$ws.Cells.Item(51, 5).Value = "1185.43126903932"
$ws.Cells.Item(52, 5).Value = "1191.95387984613"
$ws.Cells.Item(53, 5).Value = "1255.88443421343"
$ws.Cells.Item(54, 5).Value = "1301.05882786418"
$ws.Cells.Item(55, 5).Value = "1335.03049678079"
$ws.Cells.Item(56, 5).Value = "1364.68519311037"
$ws.Cells.Item(57, 5).Value = "1451.09732627394"
$ws.Cells.Item(58, 5).Value = "1501.38549955738"
$ws.Cells.Item(59, 5).Value = "1569.82984070588"
$ws.Cells.Item(60, 5).Value = "1679.39947788346"
$ws.Cells.Item(61, 5).Value = "1757.38624301723"
$ws.Cells.Item(62, 5).Value = "1831.29070046743"

# Add new rows for years 2011-2016 (rows 63-68)
$ws.Cells.Item(63, 1).Value = 800
$ws.Cells.Item(63, 2).Value = "Uganda"
$ws.Cells.Item(63, 3).Value = "GDP per Capita"
$ws.Cells.Item(63, 4).Value = 2011
$ws.Cells.Item(63, 5).Value = "1890"
$ws.Cells.Item(64, 1).Value = 800
$ws.Cells.Item(64, 2).Value = "Uganda"
$ws.Cells.Item(64, 3).Value = "GDP per Capita"
$ws.Cells.Item(64, 4).Value = 2012
$ws.Cells.Item(64, 5).Value = "1877"
$ws.Cells.Item(65, 1).Value = 800
$ws.Cells.Item(65, 2).Value = "Uganda"
$ws.Cells.Item(65, 3).Value = "GDP per Capita"
$ws.Cells.Item(65, 4).Value = 2013
$ws.Cells.Item(65, 5).Value = "1889"
$ws.Cells.Item(66, 1).Value = 800
$ws.Cells.Item(66, 2).Value = "Uganda"
$ws.Cells.Item(66, 3).Value = "GDP per Capita"
$ws.Cells.Item(66, 4).Value = 2014
$ws.Cells.Item(66, 5).Value = "1923"
$ws.Cells.Item(67, 1).Value = 800
$ws.Cells.Item(67, 2).Value = "Uganda"
$ws.Cells.Item(67, 3).Value = "GDP per Capita"
$ws.Cells.Item(67, 4).Value = 2015
$ws.Cells.Item(67, 5).Value = "1954"
$ws.Cells.Item(68, 1).Value = 800
$ws.Cells.Item(68, 2).Value = "Uganda"
$ws.Cells.Item(68, 3).Value = "GDP per Capita"
$ws.Cells.Item(68, 4).Value = 2016
$ws.Cells.Item(68, 5).Value = "1980"
